$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "stamina" column header in K1, copying the header style from
# an existing header cell (J1) so formatting (bold, border, centered) matches.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "stamina"

# Updated player data: row, Name, speed(F), agility(G), creating(H), shooting(I), stability(J), stamina(K)
$rows = @(
    @(2, "William", 51, 58, 62, 78, 86, 67),
    @(3, "Allison", 72, 65, 75, 63, 62, 69),
    @(4, "Kelli", 70, 64, 77, 70, 83, 90),
    @(5, "Duane", 78, 76, 77, 69, 80, 58),
    @(6, "Linda", 62, 71, 69, 83, 68, 85),
    @(7, "Helen", 84, 54, 73, 75, 65, 75),
    @(8, "Andrew", 56, 76, 79, 75, 65, 77),
    @(9, "Christopher", 75, 51, 89, 74, 60, 57),
    @(10, "Robert", 56, 85, 66, 76, 70, 90),
    @(11, "Adam", 77, 57, 70, 68, 82, 81),
    @(12, "David", 68, 77, 70, 68, 66, 63),
    @(13, "Denise", 70, 61, 73, 76, 69, 70),
    @(14, "Michael", 85, 76, 76, 65, 56, 78),
    @(15, "Tiffany", 72, 69, 76, 75, 81, 84),
    @(16, "Richard", 84, 71, 71, 88, 71, 89),
    @(17, "Mary", 70, 63, 72, 62, 68, 84),
    @(18, "Willis", 66, 83, 64, 65, 64, 79),
    @(19, "Kelli", 96, 45, 75, 76, 45, 59),
    @(20, "James", 63, 63, 92, 82, 64, 69),
    @(21, "James", 70, 60, 65, 58, 71, 73),
    @(22, "Annie", 55, 53, 79, 70, 69, 38),
    @(23, "Edward", 66, 71, 53, 60, 70, 72),
    @(24, "Walter", 61, 66, 73, 56, 63, 60),
    @(25, "Julie", 68, 45, 85, 90, 47, 70),
    @(26, "Olive", 64, 60, 68, 72, 78, 70),
    @(27, "Debra", 77, 58, 58, 70, 56, 69),
    @(28, "Bonnie", 70, 75, 55, 65, 88, 76),
    @(29, "Sylvia", 66, 70, 63, 72, 78, 61),
    @(30, "Brooke", 91, 72, 75, 74, 62, 74),
    @(31, "Scott", 78, 93, 89, 69, 80, 62),
    @(32, "Felicia", 72, 59, 83, 71, 87, 79),
    @(33, "Robert", 43, 69, 66, 65, 71, 77),
    @(34, "Sharon", 68, 71, 60, 82, 54, 65),
    @(35, "Muriel", 81, 81, 61, 66, 39, 71),
    @(36, "Sadie", 87, 70, 85, 76, 72, 88),
    @(37, "Robert", 68, 85, 83, 71, 78, 69),
    @(38, "Dale", 87, 66, 80, 81, 84, 72),
    @(39, "Paul", 64, 76, 61, 64, 66, 74),
    @(40, "John", 69, 80, 66, 40, 66, 80),
    @(41, "Diane", 82, 73, 56, 82, 60, 61)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]
    $ws.Cells.Item($r, 9).Value = $row[5]
    $ws.Cells.Item($r, 10).Value = $row[6]
    $ws.Cells.Item($r, 11).Value = $row[7]
}
